$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column A currently holds "79174445" as text (inline string).
# Convert it to a real number to match the target state.
$ws.Range("A7").Value = 79174445

# Append new payment row 8.
# Leading apostrophe forces text storage (phone numbers are stored as
# text here, same as the original A7 before its numeric conversion).
$ws.Range("A8").Value = "'79174445"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = "Cash"
$ws.Range("H8").Value = "2025-08-29T16:21:46"
